$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 23
$ws1.Range("F3").Value = 467
$ws1.Range("F6").Value = 10
$ws1.Range("F7").Value = 1305
$ws1.Range("F8").Value = 513
$ws1.Range("F10").Value = 1306
$ws1.Range("F11").Value = 178
$ws1.Range("F12").Value = 1092
$ws1.Range("F13").Value = 23
$ws1.Range("F16").Value = 105
$ws1.Range("F18").Value = 1662
$ws1.Range("F19").Value = 616
$ws1.Range("F20").Value = 270
$ws1.Range("F21").Value = 225
$ws1.Range("F22").Value = 2329
$ws1.Range("F23").Value = 7
$ws1.Range("F24").Value = 405
$ws1.Range("F26").Value = 928
$ws1.Range("F27").Value = 1212
$ws1.Range("F30").Value = 2823
$ws1.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202405/wvUOuzhk1715681404265.png"
$ws1.Range("F31").Value = 1619
$ws1.Range("F34").Value = 668
$ws1.Range("F36").Value = 1822
$ws1.Range("F37").Value = 892
$ws1.Range("F38").Value = 1839
$ws1.Range("F41").Value = 842
$ws1.Range("F42").Value = 40
$ws1.Range("F43").Value = 868
$ws1.Range("F44").Value = 796
$ws1.Range("F45").Value = 1015
$ws1.Range("F46").Value = 87
$ws1.Range("F49").Value = 3342

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 801
$ws2.Range("F17").Value = 9

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 23
$ws4.Range("F3").Value = 467
$ws4.Range("F8").Value = 10
$ws4.Range("F9").Value = 1305
$ws4.Range("F10").Value = 513
$ws4.Range("F12").Value = 1306
$ws4.Range("F13").Value = 178
$ws4.Range("F14").Value = 1092
$ws4.Range("F15").Value = 23
$ws4.Range("F18").Value = 105
$ws4.Range("F20").Value = 1662
$ws4.Range("F21").Value = 616
$ws4.Range("F22").Value = 270
$ws4.Range("F23").Value = 225
$ws4.Range("F24").Value = 2329
$ws4.Range("F25").Value = 405
$ws4.Range("F27").Value = 1212
$ws4.Range("F28").Value = 2823
$ws4.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202405/wvUOuzhk1715681404265.png"
$ws4.Range("F29").Value = 1619
$ws4.Range("F32").Value = 801
$ws4.Range("F34").Value = 668
$ws4.Range("F36").Value = 1822
$ws4.Range("F37").Value = 9
$ws4.Range("F38").Value = 892
$ws4.Range("F39").Value = 1839
$ws4.Range("F40").Value = 842
$ws4.Range("F41").Value = 868
$ws4.Range("F42").Value = 796
$ws4.Range("F43").Value = 1015
$ws4.Range("F44").Value = 87
$ws4.Range("F48").Value = 3342
